$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.576.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.755.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4526"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3573"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07478"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.45"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.089"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.174"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.752.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001057"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.12"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.756"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.618.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.63"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.958.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.61"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09188"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.647"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.524"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02288"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.73"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2091"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06004"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6285"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.931"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.791"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.714"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5868"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.91"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.938"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06890"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.133"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.61"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.55%  "
